$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 13): 支出, 300, 2017-09-25, 校园卡费用（76/月）
$ws.Range("B13").Value = "支出"
$ws.Range("C13").Value = 300
$ws.Range("D13").Value2 = 43003
$ws.Range("D12").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("E13").Value = "校园卡费用（76/月）"

# Update selection to match the edited cell
$ws.Range("E14").Select()
